# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto symbol list
# with the latest scraped values (GitHub Actions update run).
#
# The source cells are plain text ("inlineStr") rather than numbers/percentages
# (e.g. D2 = "265.76", E2 = "1.73%"), so a naive `.Value = "265.76"` assignment
# would let Excel's type inference turn it into a real number/percentage. To
# preserve the original text representation we briefly force a text number
# format ("@") before writing the literal string, then restore the cell's
# style to "Normal" so no residual formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "265.76"
    "E2" = "1.73%"
    "D3" = "26.83"
    "E3" = "-1.01%"
    "D4" = "4.700"
    "E4" = "-0.02%"
    "D5" = "0.06089"
    "E5" = "-1.62%"
    "D6" = "6.737"
    "E7" = "0.12%"
    "D8" = "0.9083"
    "E8" = "-0.83%"
    "D9" = "0.1407"
    "E9" = "0.19%"
    "D10" = "0.05075"
    "E10" = "8.63%"
    "D11" = "0.07096"
    "E11" = "0.12%"
    "D12" = "0.03148"
    "E12" = "-0.06%"
    "D13" = "0.09029"
    "E13" = "-0.18%"
    "D14" = "0.001527"
    "E14" = "-0.44%"
    "D15" = "0.0006069"
    "E15" = "-1.72%"
    "D16" = "0.006117"
    "E16" = "-0.20%"
    "E17" = "-0.02%"
    "D18" = "3.167"
    "E18" = "-0.31%"
    "D19" = "2.175"
    "E19" = "1.04%"
    "E21" = "-1.37%"
    "D22" = "4.091"
    "E22" = "0.27%"
    "D23" = "0.04244"
    "E23" = "0.62%"
    "E24" = "-3.16%"
    "E25" = "6.79%"
    "E26" = "0.05%"
    "E27" = "6.59%"
    "D40" = "0.03927"
    "E40" = "0.91%"
    "D41" = "0.1115"
    "E41" = "0.31%"
    "D42" = "0.004193"
    "E42" = "2.27%"
    "D43" = "0.002111"
    "E43" = "-3.33%"
    "D44" = "0.01149"
    "E44" = "-29.60%"
    "D45" = "0.00005128"
    "E45" = "-0.61%"
    "E46" = "0.06%"
    "D48" = "0.2575"
    "E48" = "54.51%"
    "E49" = "0.06%"
    "E50" = "0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
